$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Checklist" to "Session"
$ws.Name = "Session"

# Delete row 2 (180003 / Selection entry) - this shifts rows 3 and 4 up to rows 2 and 3
$ws.Rows(2).Delete()

# Update the "Type" column (E) values from "Selection" to "Scan" for the remaining data rows
$ws.Range("E2").Value = "Scan"
$ws.Range("E3").Value = "Scan"
